# Swap the "分类名称" (category name) and "单品名称" (item name) columns
# in the sumSales_year export. Column A previously held 分类名称 (always
# "食用菌") and column B held 单品名称 (the specific product name); the new
# layout puts 单品名称 in column A and 分类名称 in column B, for every row
# including the header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 1; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellB = $ws.Cells.Item($r, 2)

    $valA = $cellA.Value()
    $valB = $cellB.Value()

    $cellA.Value = $valB
    $cellB.Value = $valA
}
